# Apply cryptos list update (values refreshed by GitHub Actions run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.075.45"
$ws.Range("E2").Value = "  -0.07%  "
$ws.Range("D3").Value = "1.650.67"
$ws.Range("E3").Value = "  +0.00%  "
$ws.Range("E4").Value = "  -0.30%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "218.08"
$ws.Range("E5").Value = "  +0.36%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5209"
$ws.Range("E6").Value = "  +0.33%  "
$ws.Range("E7").Value = "  -0.31%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2634"
$ws.Range("E8").Value = "  +0.70%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06324"
$ws.Range("E9").Value = "  +0.73%  "
$ws.Range("E10").Value = "  -0.54%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07653"
$ws.Range("E11").Value = "  -1.85%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.589"
$ws.Range("E12").Value = "  +2.95%  "
$ws.Range("D13").Value = "1.634.13"
$ws.Range("E13").Value = "  -0.06%  "
$ws.Range("D14").Value = "1.877.31"
$ws.Range("E14").Value = "  +0.04%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5586"
$ws.Range("E15").Value = "  +0.77%  "
$ws.Range("D16").Value = "0.0₅8138"
$ws.Range("E16").Value = "  +1.85%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "65.20"
$ws.Range("E17").Value = "  +0.82%  "
$ws.Range("D18").Value = "26.050.38"
$ws.Range("E18").Value = "  -0.13%  "
$ws.Range("E19").Value = "  -0.34%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.622"
$ws.Range("E20").Value = "  -0.11%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.47"
$ws.Range("E21").Value = "  +3.99%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "190.97"
$ws.Range("E22").Value = "  -1.74%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.905"
$ws.Range("E23").Value = "  -0.64%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.003"
$ws.Range("E24").Value = "  -0.34%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.13"
$ws.Range("E25").Value = "  -1.80%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1186"
$ws.Range("E26").Value = "  -1.53%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.198"
$ws.Range("E27").Value = "  +0.27%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.88"
$ws.Range("E28").Value = "  +0.12%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.512"
$ws.Range("E29").Value = "  +2.45%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05430"
$ws.Range("E30").Value = "  -3.03%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.265"
$ws.Range("E31").Value = "  +0.13%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.438"
$ws.Range("E32").Value = "  -1.05%  "
$ws.Range("E33").Value = "  -1.06%  "
$ws.Range("E34").Value = "  -2.28%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.425"
$ws.Range("E35").Value = "  +0.88%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.782"
$ws.Range("E36").Value = "  -0.63%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9447"
$ws.Range("E37").Value = "  -0.37%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5628"
$ws.Range("E38").Value = "  -0.32%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01581"
$ws.Range("E39").Value = "  +0.25%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.853"
$ws.Range("E40").Value = "  -2.01%  "
$ws.Range("E41").Value = "  -0.23%  "
$ws.Range("D42").Value = "1.027.82"
$ws.Range("E42").Value = "  -3.05%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8226"
$ws.Range("E43").Value = "  -1.92%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.05"
$ws.Range("E44").Value = "  -1.49%  "
$ws.Range("D45").Value = "1.786.63"
$ws.Range("E45").Value = "  -0.09%  "
$ws.Range("D46").Value = "0.0₈111"
$ws.Range("E46").Value = "  +3.98%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "57.28"
$ws.Range("E47").Value = "  +0.39%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.001"
$ws.Range("E48").Value = "  -0.51%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4330"
$ws.Range("E49").Value = "  -0.08%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.952"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05144"
$ws.Range("E51").Value = "  -3.81%  "
